$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 223, pushing existing rows 223:355 down to 224:356.
$ws.Rows.Item(223).Insert()

# Populate the newly inserted row 223 with the new data point.
$ws.Cells.Item(223, 1).Value = 10
$ws.Cells.Item(223, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(223, 3).Value = "La Araucanía"
$ws.Cells.Item(223, 4).Value = 45176
$ws.Cells.Item(223, 5).Value = 9
$ws.Cells.Item(223, 6).Value = 100114007
$ws.Cells.Item(223, 7).Value = "Jengibre"
$ws.Cells.Item(223, 8).Value = "Sin especificar"
$ws.Cells.Item(223, 9).Value = "Primera"
$ws.Cells.Item(223, 10).Value = 200
$ws.Cells.Item(223, 11).Value = 22000
$ws.Cells.Item(223, 12).Value = 24000
$ws.Cells.Item(223, 13).Value = 23200
$ws.Cells.Item(223, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(223, 15).Value = "Perú"
$ws.Cells.Item(223, 16).Value = 1785
$ws.Cells.Item(223, 17).Value = 13
$ws.Cells.Item(223, 18).Value = "Hortaliza"
